$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 656
$ws.Range("I2").Value = 1790
$ws.Range("J2").Value = 7224
$ws.Range("K2").Value = 46
$ws.Range("L2").Value = 1932
$ws.Range("M2").Value = 149
$ws.Range("N2").Value = 1291
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 25
$ws.Range("Q2").Value = 14
$ws.Range("R2").Value = 87
$ws.Range("S2").Value = 792
$ws.Range("T2").Value = 1302
$ws.Range("U2").Value = 94
$ws.Range("V2").Value = 11355
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 11230
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 152
$ws.Range("AA2").Value = 73
